# Rename "Device" sheet to "Apparatus" and update its labels/header text
# to match (form, simulink model generation, function naming all refer to
# "Device" -> "Apparatus").

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Device")
$ws.Name = "Apparatus"

# Sheet summary / header text
$ws.Range("A1").Value = "This sheet summarizes the apparatuses connected to buses."
$ws.Range("B2").Value = "Apparatus type"
$ws.Range("C2").Value = "Apparatus parameters"

# Bring the renamed sheet into focus, as the author had it selected when
# the workbook was last saved.
$ws.Activate()
